# Example2.xlsx manual-rewrite edit
# The author replaced several "departure time" formulas (=C# + <hours>) in the
# Aircraft_scheduling sheet with the literal time-of-day values that those
# formulas used to compute (minus the inadvertent extra day the old formula
# carried when the duration pushed the result past midnight). The Turn
# Around Time column (F) keeps its shared MOD(E-C,1) formula and Excel
# recalculates it automatically once E no longer has a formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Aircraft_scheduling")

$ws.Range("E6").Value  = 0.73932870370370374
$ws.Range("E11").Value = 0.35043981481481484
$ws.Range("E12").Value = 0.17569444444444446
$ws.Range("E13").Value = 0.41319444444444442
$ws.Range("E14").Value = 0.53758101851851847
$ws.Range("E15").Value = 0.36402777777777778
$ws.Range("E16").Value = 0.48837962962962966
$ws.Range("E17").Value = 0.3972222222222222
$ws.Range("E18").Value = 0.66629629629629628

# Restore the recorded selection/active cell for the sheet.
$ws.Range("E19").Select()
